# Update the "Förändrad" (Changed) date column for rows 2-6 from
# 2023-11-13 (serial 45243) to 2023-11-14 (serial 45244).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C6").Value = 45244
